# Updates cryptos list data (prices + volume%) per upstream GitHub Actions run.
# Most of column D contains numeric-looking strings (e.g. "1.00", "0.998") that
# must be preserved as literal text, matching the existing inlineStr cells in
# the workbook. We force text via NumberFormat "@" before assignment, then
# ClearFormats() to drop the residual style so the cell's style index is left
# unchanged (back to the sheet's default, unstyled look).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$ws.Range("D2").Value = '67.449.92'
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").Value = '3.286.85'
$ws.Range("E3").Value = '  +0.05%  '
Set-TextValue "D4" '0.998'
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue "D5" '574.48'
$ws.Range("E5").Value = '  -0.43%  '
Set-TextValue "D6" '177.09'
$ws.Range("E6").Value = '  -2.58%  '
$ws.Range("E7").Value = '  +0.09%  '
Set-TextValue "D8" '0.584'
$ws.Range("E8").Value = '  +3.00%  '
$ws.Range("D9").Value = '3.279.52'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").Value = '  -0.05%  '
Set-TextValue "D11" '0.573'
$ws.Range("E11").Value = '  +1.21%  '
Set-TextValue "D12" '45.50'
$ws.Range("E12").Value = '  -1.59%  '
Set-TextValue "D13" '0.0000269'
$ws.Range("E13").Value = '  +2.51%  '
Set-TextValue "D14" '702.75'
$ws.Range("E14").Value = '  +12.14%  '
$ws.Range("D15").Value = '3.812.67'
$ws.Range("E15").Value = '  +0.27%  '
Set-TextValue "D16" '8.34'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '67.539.26'
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").Value = '3.282.81'
$ws.Range("E19").Value = '  -0.01%  '
Set-TextValue "D20" '17.38'
$ws.Range("E20").Value = '  -1.15%  '
Set-TextValue "D21" '10.77'
$ws.Range("E21").Value = '  -0.77%  '
Set-TextValue "D22" '0.890'
$ws.Range("E22").Value = '  +0.92%  '
Set-TextValue "D23" '16.94'
$ws.Range("E23").Value = '  -5.44%  '
Set-TextValue "D24" '5.14'
$ws.Range("E24").Value = '  +3.92%  '
Set-TextValue "D25" '98.64'
$ws.Range("E25").Value = '  -0.68%  '
$ws.Range("E26").Value = '  -0.82%  '
Set-TextValue "D27" '2.72'
$ws.Range("E27").Value = '  +0.34%  '
Set-TextValue "D28" '9.31'
$ws.Range("E28").Value = '  +0.09%  '
Set-TextValue "D29" '32.99'
$ws.Range("E29").Value = '  +7.95%  '
Set-TextValue "D30" '8.43'
$ws.Range("E30").Value = '  +1.32%  '
Set-TextValue "D31" '6.67'
$ws.Range("E31").Value = '  +3.69%  '
Set-TextValue "D32" '580.46'
$ws.Range("E32").Value = '  +4.40%  '
$ws.Range("D33").Value = '3.878.52'
$ws.Range("E33").Value = '  +2.03%  '
Set-TextValue "D34" '10.80'
$ws.Range("E34").Value = '  +0.29%  '
$ws.Range("E35").Value = '  +0.45%  '
Set-TextValue "D36" '1.00'
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D37" '3.33'
$ws.Range("E37").Value = '  -8.25%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D38" '55.30'
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("E40").Value = '  +1.80%  '
Set-TextValue "D41" '3.13'
$ws.Range("E41").Value = '  +0.15%  '
Set-TextValue "D42" '31.85'
$ws.Range("E42").Value = '  -1.16%  '
Set-TextValue "D43" '3.34'
$ws.Range("E43").Value = '  -2.24%  '
$ws.Range("D44").Value = '0.0₃0674'
$ws.Range("E44").Value = '  +0.45%  '
Set-TextValue "D45" '0.329'
$ws.Range("E45").Value = '  +0.28%  '
Set-TextValue "D46" '0.0407'
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("E47").Value = '  +1.75%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D48" '1.01'
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D49" '1.39'
$ws.Range("E49").Value = '  +10.53%  '
Set-TextValue "D50" '2.52'
$ws.Range("E50").Value = '  +1.26%  '
Set-TextValue "D51" '127.45'
$ws.Range("E51").Value = '  -0.53%  '
